$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 482; $row++) {
    $ws.Cells.Item($row, 3).Value = 45205
}
